$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 115.55556
$ws.Range("I6").Value = 115.55556
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 346.66668
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -234.66668
$ws.Range("N6").ClearContents()
$ws.Range("H12").Value = 840
$ws.Range("J12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("N12").Value = -2840
$ws.Range("H19").Value = 1285.909
$ws.Range("J19").Value = 1995.8334
$ws.Range("L19").Value = 1995.8334
$ws.Range("N19").Value = -2345.8334
$ws.Range("H28").Value = 914.8
$ws.Range("I28").Value = 914.8
$ws.Range("K28").Value = 914.8
$ws.Range("M28").Value = -429.8
$ws.Range("H31").Value = 1360.4
$ws.Range("I31").Value = 1360.4
$ws.Range("K31").Value = 4081.2
$ws.Range("M31").Value = -3851.2
$ws.Range("H33").Value = 948.0909
$ws.Range("I33").Value = 1005
$ws.Range("K33").Value = 1005
$ws.Range("M33").Value = -776
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 6000
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -5881
$ws.Range("N46").Value = -6238
$ws.Range("H52").Value = 7418.3
$ws.Range("I52").Value = 592.5
$ws.Range("J52").Value = 9124.75
$ws.Range("K52").Value = 1777.5
$ws.Range("L52").Value = 27374.25
$ws.Range("M52").Value = -1617.5
$ws.Range("N52").Value = -27694.25
$ws.Range("H60").Value = 2000
$ws.Range("I60").Value = 2000
$ws.Range("J60").Value = 2000
$ws.Range("K60").Value = 6000
$ws.Range("L60").Value = 6000
$ws.Range("M60").Value = -5516
$ws.Range("N60").Value = -6968
$ws.Range("I62").Value = 5152.5454
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 5152.5454
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -4528.5454
$ws.Range("N62").Value = -4148
$ws.Range("I65").Value = 5152.5454
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 25762.727
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -22642.727
$ws.Range("N65").Value = -20740
$ws.Range("H70").Value = 5191
$ws.Range("I70").Value = 4951
$ws.Range("J70").Value = 5311
$ws.Range("K70").Value = 14853
$ws.Range("L70").Value = 15933
$ws.Range("M70").Value = -14583
$ws.Range("N70").Value = -16473
$ws.Range("H73").Value = 5191
$ws.Range("I73").Value = 4951
$ws.Range("J73").Value = 5311
$ws.Range("K73").Value = 14853
$ws.Range("L73").Value = 15933
$ws.Range("M73").Value = -13917
$ws.Range("N73").Value = -17805
$ws.Range("H80").Value = 3472.8462
$ws.Range("I80").Value = 1125.9231
$ws.Range("J80").Value = 4646.3076
$ws.Range("K80").Value = 3377.7693
$ws.Range("L80").Value = 13938.9228
$ws.Range("M80").Value = -2379.7693
$ws.Range("N80").Value = -15934.9228
$ws.Range("H83").Value = 3472.8462
$ws.Range("I83").Value = 1125.9231
$ws.Range("J83").Value = 4646.3076
$ws.Range("K83").Value = 10133.3079
$ws.Range("L83").Value = 41816.7684
$ws.Range("M83").Value = -5141.3079
$ws.Range("N83").Value = -51800.7684
$ws.Range("H98").Value = 3093.4583
$ws.Range("I98").Value = 2913.65
$ws.Range("K98").Value = 2913.65
$ws.Range("M98").Value = -1415.65
$ws.Range("H106").Value = 166668290
$ws.Range("I106").Value = 200001550
$ws.Range("K106").Value = 200001550
$ws.Range("M106").Value = -200000919
$ws.Range("H111").Value = 5275.7144
$ws.Range("I111").Value = 5424.5
$ws.Range("K111").Value = 16273.5
$ws.Range("M111").Value = -13206.5
$ws.Range("H116").Value = 5824.3335
$ws.Range("J116").Value = 5874.3335
$ws.Range("L116").Value = 5874.3335
$ws.Range("N116").Value = -12758.3335
$ws.Range("H122").Value = 3093.4583
$ws.Range("I122").Value = 2913.65
$ws.Range("K122").Value = 8740.950000000001
$ws.Range("M122").Value = -6290.950000000001
$ws.Range("H132").Value = 4157.5435
$ws.Range("I132").Value = 3906.5898
$ws.Range("J132").Value = 5555.7144
$ws.Range("K132").Value = 11719.7694
$ws.Range("L132").Value = 16667.1432
$ws.Range("M132").Value = -9189.769400000001
$ws.Range("N132").Value = -21727.1432
$ws.Range("H135").Value = 1758.0526
$ws.Range("I135").Value = 1356.5
$ws.Range("J135").Value = 3899.6667
$ws.Range("K135").Value = 12208.5
$ws.Range("L135").Value = 35097.0003
$ws.Range("M135").Value = -9673.5
$ws.Range("N135").Value = -40167.0003
$ws.Range("H136").Value = 85493.336
$ws.Range("J136").Value = 85493.336
$ws.Range("L136").Value = 85493.336
$ws.Range("N136").Value = -95693.336
$ws.Range("H137").Value = 78757.17
$ws.Range("I137").Value = 78757.17
$ws.Range("K137").Value = 236271.51
$ws.Range("M137").Value = -233721.51

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H61").Value = 3738.0938
$ws.Range("I61").Value = 3566.1333
$ws.Range("J61").Value = 3889.8235
$ws.Range("K61").Value = 3566.1333
$ws.Range("L61").Value = 3889.8235
$ws.Range("M61").Value = -3354.1333
$ws.Range("N61").Value = -4313.8235
$ws.Range("H74").Value = 195626.36
$ws.Range("I74").Value = 114218.125
$ws.Range("K74").Value = 114218.125
$ws.Range("M74").Value = -113344.125
$ws.Range("H77").Value = 195626.36
$ws.Range("I77").Value = 114218.125
$ws.Range("K77").Value = 571090.625
$ws.Range("M77").Value = -566722.625
$ws.Range("H102").Value = 3257.9333
$ws.Range("I102").Value = 1483.6842
$ws.Range("J102").Value = 6322.5454
$ws.Range("K102").Value = 1483.6842
$ws.Range("L102").Value = 6322.5454
$ws.Range("M102").Value = 138.3158000000001
$ws.Range("N102").Value = -9566.545399999999
$ws.Range("H110").Value = 1417.0435
$ws.Range("I110").Value = 1260.7778
$ws.Range("K110").Value = 1260.7778
$ws.Range("M110").Value = 784.2221999999999
$ws.Range("H132").Value = 3752.5833
$ws.Range("I132").Value = 3238.3076
$ws.Range("J132").Value = 4043.261
$ws.Range("K132").Value = 9714.9228
$ws.Range("L132").Value = 12129.783
$ws.Range("M132").Value = -7184.9228
$ws.Range("N132").Value = -17189.783
$ws.Range("H136").Value = 3738.0938
$ws.Range("I136").Value = 3566.1333
$ws.Range("J136").Value = 3889.8235
$ws.Range("K136").Value = 10698.3999
$ws.Range("L136").Value = 11669.4705
$ws.Range("M136").Value = -8148.3999
$ws.Range("N136").Value = -16769.4705
$ws.Range("H141").Value = 96199.2
$ws.Range("J141").Value = 72749.25
$ws.Range("L141").Value = 72749.25
$ws.Range("N141").Value = -83109.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7437.2607
$ws.Range("I86").Value = 7907.85
$ws.Range("K86").Value = 7907.85
$ws.Range("M86").Value = -6784.85
$ws.Range("H89").Value = 7437.2607
$ws.Range("I89").Value = 7907.85
$ws.Range("K89").Value = 39539.25
$ws.Range("M89").Value = -33923.25
$ws.Range("H94").Value = 7080.032
$ws.Range("I94").Value = 1052.8928
$ws.Range("J94").Value = 63333.332
$ws.Range("K94").Value = 1052.8928
$ws.Range("L94").Value = 63333.332
$ws.Range("M94").Value = -601.8928000000001
$ws.Range("N94").Value = -64235.332
$ws.Range("H105").Value = 5064.5557
$ws.Range("I105").Value = 5321.25
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 5321.25
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -3574.25
$ws.Range("N105").Value = -6505
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H134").Value = 4457.735
$ws.Range("J134").Value = 6253.9644
$ws.Range("L134").Value = 18761.8932
$ws.Range("N134").Value = -23831.8932
$ws.Range("H137").Value = 72499.836
$ws.Range("J137").Value = 72499.836
$ws.Range("L137").Value = 72499.836
$ws.Range("N137").Value = -82699.836

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1688
$ws.Range("I6").Value = 1711.5
$ws.Range("K6").Value = 1711.5
$ws.Range("M6").Value = -1598.5
$ws.Range("H7").Value = 259.73334
$ws.Range("I7").Value = 242.57143
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 242.57143
$ws.Range("L7").Value = 500
$ws.Range("M7").Value = -129.57143
$ws.Range("N7").Value = -726
$ws.Range("H16").Value = 2612.5
$ws.Range("I16").Value = 1793.5714
$ws.Range("J16").Value = 3759
$ws.Range("K16").Value = 1793.5714
$ws.Range("L16").Value = 3759
$ws.Range("M16").Value = -1506.5714
$ws.Range("N16").Value = -4333
$ws.Range("H41").Value = 8186.2856
$ws.Range("I41").Value = 3782.6365
$ws.Range("J41").Value = 24333
$ws.Range("K41").Value = 3782.6365
$ws.Range("L41").Value = 24333
$ws.Range("M41").Value = -3354.6365
$ws.Range("N41").Value = -25189
$ws.Range("H58").Value = 2057
$ws.Range("I58").Value = 1739
$ws.Range("J58").Value = 3965
$ws.Range("K58").Value = 1739
$ws.Range("L58").Value = 3965
$ws.Range("M58").Value = -1536
$ws.Range("N58").Value = -4371
$ws.Range("H59").Value = 43166.555
$ws.Range("J59").Value = 42312.375
$ws.Range("L59").Value = 42312.375
$ws.Range("N59").Value = -44602.375
$ws.Range("H68").Value = 46666.668
$ws.Range("I68").Value = 40000
$ws.Range("K68").Value = 40000
$ws.Range("M68").Value = -39251
$ws.Range("H71").Value = 46666.668
$ws.Range("I71").Value = 40000
$ws.Range("K71").Value = 120000
$ws.Range("M71").Value = -116256
$ws.Range("H94").Value = 1481.7778
$ws.Range("J94").Value = 1504.5
$ws.Range("L94").Value = 1504.5
$ws.Range("N94").Value = -2406.5
$ws.Range("H97").Value = 59890
$ws.Range("J97").Value = 59890
$ws.Range("L97").Value = 59890
$ws.Range("N97").Value = -61872
$ws.Range("H107").Value = 27028386
$ws.Range("I107").Value = 1351.3214
$ws.Range("K107").Value = 1351.3214
$ws.Range("M107").Value = 568.6786
$ws.Range("H113").Value = 2612.5
$ws.Range("I113").Value = 1793.5714
$ws.Range("J113").Value = 3759
$ws.Range("K113").Value = 1793.5714
$ws.Range("L113").Value = 3759
$ws.Range("M113").Value = 376.4286
$ws.Range("N113").Value = -8099
$ws.Range("H132").Value = 25426.46
$ws.Range("I132").Value = 1540.6154
$ws.Range("K132").Value = 4621.8462
$ws.Range("M132").Value = -2091.8462
$ws.Range("H134").Value = 30807.031
$ws.Range("J134").Value = 5498
$ws.Range("L134").Value = 16494
$ws.Range("N134").Value = -21564
$ws.Range("H136").Value = 2057
$ws.Range("I136").Value = 1739
$ws.Range("J136").Value = 3965
$ws.Range("K136").Value = 5217
$ws.Range("L136").Value = 11895
$ws.Range("M136").Value = -2667
$ws.Range("N136").Value = -16995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 155.82353
$ws.Range("J26").Value = 124.166664
$ws.Range("L26").Value = 372.499992
$ws.Range("N26").Value = -948.499992
$ws.Range("H107").Value = 566.4231
$ws.Range("J107").Value = 758
$ws.Range("L107").Value = 2274
$ws.Range("N107").Value = -6114
$ws.Range("H122").Value = 1080.0555
$ws.Range("I122").Value = 1014.63635
$ws.Range("J122").Value = 1182.8572
$ws.Range("K122").Value = 9131.727150000001
$ws.Range("L122").Value = 10645.7148
$ws.Range("M122").Value = -6681.727150000001
$ws.Range("N122").Value = -15545.7148
$ws.Range("H131").Value = 7719268
$ws.Range("I131").Value = 4387486.5
$ws.Range("J131").Value = 9527950
$ws.Range("K131").Value = 13162459.5
$ws.Range("L131").Value = 28583850
$ws.Range("M131").Value = -13157419.5
$ws.Range("N131").Value = -28593930
$ws.Range("H132").Value = 2672.1892
$ws.Range("I132").Value = 1092.091
$ws.Range("K132").Value = 9828.819
$ws.Range("M132").Value = -7298.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9079.714
$ws.Range("I70").Value = 11538.077
$ws.Range("J70").Value = 5084.875
$ws.Range("K70").Value = 11538.077
$ws.Range("L70").Value = 5084.875
$ws.Range("M70").Value = -11268.077
$ws.Range("N70").Value = -5624.875
$ws.Range("H73").Value = 9079.714
$ws.Range("I73").Value = 11538.077
$ws.Range("J73").Value = 5084.875
$ws.Range("K73").Value = 11538.077
$ws.Range("L73").Value = 5084.875
$ws.Range("M73").Value = -10602.077
$ws.Range("N73").Value = -6956.875
$ws.Range("H80").Value = 33469462
$ws.Range("I80").Value = 55558320
$ws.Range("J80").Value = 336178.5
$ws.Range("K80").Value = 55558320
$ws.Range("L80").Value = 336178.5
$ws.Range("M80").Value = -55557322
$ws.Range("N80").Value = -338174.5
$ws.Range("H83").Value = 33469462
$ws.Range("I83").Value = 55558320
$ws.Range("J83").Value = 336178.5
$ws.Range("K83").Value = 277791600
$ws.Range("L83").Value = 1680892.5
$ws.Range("M83").Value = -277786608
$ws.Range("N83").Value = -1690876.5
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492
$ws.Range("H107").Value = 1314.9474
$ws.Range("I107").Value = 1525.9231
$ws.Range("K107").Value = 1525.9231
$ws.Range("M107").Value = 394.0769
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 50000
$ws.Range("L110").Value = 50000
$ws.Range("N110").Value = -58180
$ws.Range("H113").Value = 3680.5
$ws.Range("I113").Value = 3680.5
$ws.Range("K113").Value = 3680.5
$ws.Range("M113").Value = -1510.5
$ws.Range("H116").Value = 120321
$ws.Range("J116").Value = 120321
$ws.Range("L116").Value = 120321
$ws.Range("N116").Value = -129499
$ws.Range("H132").Value = 2923.4211
$ws.Range("I132").Value = 2803.2354
$ws.Range("J132").Value = 3945
$ws.Range("K132").Value = 8409.706200000001
$ws.Range("L132").Value = 11835
$ws.Range("M132").Value = -5879.706200000001
$ws.Range("N132").Value = -16895
$ws.Range("H136").Value = 41311
$ws.Range("J136").Value = 41311
$ws.Range("L136").Value = 123933
$ws.Range("N136").Value = -129033
$ws.Range("H139").Value = 96333.336
$ws.Range("J139").Value = 96333.336
$ws.Range("L139").Value = 96333.336
$ws.Range("N139").Value = -106613.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 112340.25
$ws.Range("I22").Value = 112340.25
$ws.Range("K22").Value = 112340.25
$ws.Range("M22").Value = -112045.25
$ws.Range("H27").Value = 112340.25
$ws.Range("I27").Value = 112340.25
$ws.Range("K27").Value = 112340.25
$ws.Range("M27").Value = -112233.25
$ws.Range("H116").Value = 72340
$ws.Range("J116").Value = 72340
$ws.Range("L116").Value = 72340
$ws.Range("N116").Value = -81518
$ws.Range("H119").Value = 102500
$ws.Range("J119").Value = 102500
$ws.Range("L119").Value = 102500
$ws.Range("N119").Value = -112176
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 6756.3076
$ws.Range("I122").Value = 4410
$ws.Range("J122").Value = 8222.75
$ws.Range("K122").Value = 13230
$ws.Range("L122").Value = 24668.25
$ws.Range("M122").Value = -10780
$ws.Range("N122").Value = -29568.25
$ws.Range("H132").Value = 5722.0425
$ws.Range("I132").Value = 5219.815
$ws.Range("J132").Value = 7317.353
$ws.Range("K132").Value = 15659.445
$ws.Range("L132").Value = 21952.059
$ws.Range("M132").Value = -13129.445
$ws.Range("N132").Value = -27012.059
$ws.Range("H134").Value = 101096.71
$ws.Range("J134").Value = 101096.71
$ws.Range("L134").Value = 101096.71
$ws.Range("N134").Value = -111236.71

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 1737.5
$ws.Range("I40").Value = 1737.5
$ws.Range("K40").Value = 1737.5
$ws.Range("M40").Value = -1588.5
$ws.Range("H70").Value = 31665.834
$ws.Range("J70").Value = 39998.332
$ws.Range("L70").Value = 39998.332
$ws.Range("N70").Value = -40628.332
$ws.Range("H73").Value = 31665.834
$ws.Range("J73").Value = 39998.332
$ws.Range("L73").Value = 39998.332
$ws.Range("N73").Value = -42182.332
$ws.Range("H107").Value = 1666.6666
$ws.Range("I107").Value = 1500
$ws.Range("K107").Value = 4500
$ws.Range("M107").Value = -2580
$ws.Range("H109").Value = 44027
$ws.Range("J109").Value = 44027
$ws.Range("L109").Value = 44027
$ws.Range("N109").Value = -46801
$ws.Range("H117").Value = 64090
$ws.Range("J117").Value = 64090
$ws.Range("L117").Value = 64090
$ws.Range("N117").Value = -73268
$ws.Range("H122").Value = 3262.4783
$ws.Range("I122").Value = 2598.0588
$ws.Range("K122").Value = 7794.176399999999
$ws.Range("M122").Value = -5344.176399999999
$ws.Range("H132").Value = 16675.938
$ws.Range("I132").Value = 3220.0962
$ws.Range("J132").Value = 70499.30499999999
$ws.Range("K132").Value = 9660.2886
$ws.Range("L132").Value = 211497.915
$ws.Range("M132").Value = -7130.2886
$ws.Range("N132").Value = -216557.915
$ws.Range("H137").Value = 64715
$ws.Range("J137").Value = 64715
$ws.Range("L137").Value = 64715
$ws.Range("N137").Value = -74915
